# "Generate Report for Archive"
#
# The localization status report is regenerated: every cell that held the
# old "Ready for handoff" status string now reads "In Translation" (the
# Overview sheet's per-language roll-up columns, plus the Status column on
# each per-language detail sheet). Because the new text is shorter, the
# Status-related columns that were sized to fit the old text are re-sized
# to fit the new, narrower text.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
# Target rendered column width (post-edit) is ~13.41 chars; ColumnWidth is a
# pre-padding character width that this host then quantizes/pads, so feed it
# the value that lands on the nearest attainable stored width.
$newWidth = 12.5

# --- Overview sheet: zh-cn (col E) / de-de (col F) status roll-up ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = $newWidth
$overview.Columns.Item(6).ColumnWidth = $newWidth

# --- zh-cn sheet: Status column (col C) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Columns.Item(3).ColumnWidth = $newWidth

# --- de-de sheet: Status column (col C) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Columns.Item(3).ColumnWidth = $newWidth
